# Ajout de plusieurs liste de questions
# Renumber the question list (1 -> 5, 2 -> 6, 3 -> 7, and related sub-labels)

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Question 1 -> 5
$ws.Range("A2").Value = 5
$ws.Range("A3").Value = "5a"
$ws.Range("A5").Value = "5b"
$ws.Range("B5").Value = "5;q5a"

# Question 2 -> 6
$ws.Range("A8").Value = 6
$ws.Range("B10").Value = 6
# Force "6.1" to be stored as text (not auto-converted to a number)
$ws.Range("A10").NumberFormat = "@"
$ws.Range("A10").Value = "6.1"
$ws.Range("A10").NumberFormat = "General"

# Question 3 -> 7
$ws.Range("A12").Value = 7

# Update the active selection to B3 (as reflected in the saved sheet view)
$ws.Range("B3").Select()
